$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Ali%Dadgari%NULL%1,                            Tengku%Aizan Hamid%NULL%1,                            Mohammad Nazrul%Hakim%NULL%1,                            Reza%Chaman%NULL%1,                            Seyed Abbas%Mousavi%NULL%1,                            Lim%Poh Hin%NULL%1,                            Leila%Dadvar%NULL%1]'
$ws.Range("E3").Value = '[ Christopher J.%Dondzila%null%1,                  Ann M.%Swartz%null%1,                  Kevin G.%Keenan%null%1,                  Amy E.%Harley%null%1,                  Razia%Azen%null%1,                  Scott J.%Strath%null%1]'
$ws.Range("E4").Value = '[Ryoichi%Ema%NULL%1,                            Shunsuke%Ohki%NULL%1,                            Hirokazu%Takayama%NULL%1,                            Yuji%Kobayashi%NULL%1,                            Ryota%Akagi%NULL%1]'
$ws.Range("E5").Value = '[Tsung-Jen%Hsieh%tjhsieh@nhri.edu.tw%1,                            Shin-Chang%Su%maxmax0220@hotmail.com%1,                            Chun-Wei%Chen%chunwei.peter@gmail.com%1,                            Yaw-Wen%Kang%castle.kang@gmail.com%1,                            Ming-Hsia%Hu%mhh@ntu.edu.tw%1,                            Li-Lin%Hsu%llh@firdi.org.tw%1,                            Szu-Yun%Wu%s.wu@ibms.sinica.edu.tw%1,                            Likwang%Chen%likwang@nhri.edu.tw%1,                            Hsing-Yi%Chang%hsingyi@nhri.edu.tw%1,                            Shao-Yuan%Chuang%chuangsy@nhri.edu.tw%1,                            Wen-Harn%Pan%pan@ibms.sinica.edu.tw%1,                            Chih-Cheng%Hsu%cch@nhri.edu.tw%1]'
$ws.Range("E6").Value = '[Steve%Iliffe%xref no email%1,       Denise%Kendrick%xref no email%1,       Richard%Morris%xref no email%1,       Mark%Griffin%xref no email%1,       Deborah%Haworth%xref no email%1,       Hannah%Carpenter%xref no email%1,       Tahir%Masud%xref no email%1,       Dawn A%Skelton%xref no email%1,       Susie%Dinan-Young%xref no email%1,       Ann%Bowling%xref no email%1,       Heather%Gage%xref no email%1]'
$ws.Range("E7").Value = '[Nicole%Kahle%xref no email%1,       Michael A.%Tevald%xref no email%1]'
$ws.Range("E9").Value = '[Andr\u00e9%Lacroix%xref no email%1,       Reto W.%Kressig%xref no email%1,       Thomas%Muehlbauer%xref no email%1,       Yves J.%Gschwind%xref no email%1,       Barbara%Pfenninger%xref no email%1,       Othmar%Bruegger%xref no email%1,       Urs%Granacher%xref no email%1]'
$ws.Range("E10").Value = '[Ahamed%Yasmin%coreGivesNoEmail%1,                          Close%Jacqueline%coreGivesNoEmail%1,                          Cook%Wendy L.%coreGivesNoEmail%1,                          Donaldson%Meghan G.%coreGivesNoEmail%1,                          Graf%Peter%coreGivesNoEmail%1,                          Khan%Karim M.%coreGivesNoEmail%1,                          Liu-Ambrose%Teresa%coreGivesNoEmail%1,                          Lord%Stephen R.%coreGivesNoEmail%1]'
$ws.Range("E11").Value = '[Kohei%Maruya%NULL%1,                            Yasuyoshi%Asakawa%NULL%1,                            Hideaki%Ishibashi%NULL%1,                            Hiroaki%Fujita%NULL%1,                            Tomoyuki%Arai%NULL%1,                            Haruyasu%Yamaguchi%NULL%1]'
$ws.Range("E12").Value = '[M. E.%Nelson%xref no email%1,       J. E.%Layne%xref no email%1,       M. J.%Bernstein%xref no email%1,       A.%Nuernberger%xref no email%1,       C.%Castaneda%xref no email%1,       D.%Kaliton%xref no email%1,       J.%Hausdorff%xref no email%1,       J. O.%Judge%xref no email%1,       D. M.%Buchner%xref no email%1,       R.%Roubenoff%xref no email%1,       M. A.%Fiatarone Singh%xref no email%1]'
$ws.Range("E14").Value = '[Oliver J.%Perkin%o.j.perkin@bath.ac.uk%1,                            Polly M.%McGuigan%NULL%2,                            Polly M.%McGuigan%NULL%0,                            Keith A.%Stokes%NULL%1]'
$ws.Range("E15").Value = '[Maria%Tsekoura%NULL%1,                            Evdokia%Billis%NULL%1,                            Elias%Tsepis%NULL%1,                            Zacharias%Dimitriadis%NULL%1,                            Charalampos%Matzaroglou%NULL%1,                            Minos%Tyllianakis%NULL%1,                            Elias%Panagiotopoulos%NULL%1,                            John%Gliatis%NULL%1]'
$ws.Range("E16").Value = '[ Sonja%Vestergaard%null%1,                  Christian%Kronborg%null%1,                  Lis%Puggaard%null%1]'
